# "nova pagina so parecer final" - add a page break near the end of the
# sheet so a new printed page only shows the final-looking rows, and
# tweak the page layout (row 4 height, top margin) plus scroll the view
# down to where the new page break lives.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("con")
$ws.Activate()

# --- Row 4 (the big merged header row) is slightly shorter now ---
$ws.Rows.Item(4).RowHeight = 177.75

# --- Give the sheet a bit more room at the top of the page ---
$ws.PageSetup.TopMargin = $excel.Application.InchesToPoints(0.39370078740157483)

# --- Manual horizontal page breaks every 21 rows from row 25 down to
#     109, so each printed page ends right after a 3-row data block ---
$breakBeforeRows = @(26, 47, 68, 89, 110)
foreach ($r in $breakBeforeRows) {
    $ws.HPageBreaks.Add($ws.Rows.Item($r)) | Out-Null
}

# --- Scroll the frozen view down near the new page break and select
#     the row that marks the final page ---
$win = $excel.ActiveWindow
$win.ScrollRow = 115
$win.ScrollColumn = 1
$ws.Rows.Item(110).Select()
